$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B14" = 1052
    "C14" = 3633
    "D14" = 14604
    "E14" = 3175
    "F14" = 4175
    "G14" = 11244
    "H14" = 60508
    "I14" = 6487
    "J14" = 4798
    "K14" = 10339
    "L14" = 3764
    "M14" = 1883
    "N14" = 4718
    "O14" = 854
    "P14" = 1597
    "Q14" = 133018
    "R14" = 78
    "S14" = 12636
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
